# Update to v1.8: bump game/table version numbers and append a new
# update-log entry, mirroring the existing "m0.1 / 完成..." pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Game version 1.4 -> 1.5, table version 0.1 -> 0.2
$ws.Range("G3").Value = 1.5
$ws.Range("H3").Value = 0.2

# New changelog rows, following the existing style used by A14/A15
$ws.Range("A16").Value = "m0.2"
$ws.Range("A17").Value = "更新游戏1.5内容"

# Move the active selection the way the author left it
$ws.Range("H4").Select()
